$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused rows 8-10 (table shrinks from 9 data rows to 6)
$ws.Range("A8:T10").EntireRow.Delete()

# New TPM-derived values for rows 2-7.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the numeric statistics columns.

$data = @(
  @{ Row=2;  A="FAPs";  B="Angpt1"; C="Tek"; D="ECs";
     E=3; F=1; G=12.95459633333333; H=38.863789; I=0.8906505749177925; J=0.8906505749177924;
     K=3; L=1; M=53.457377; N=160.372131; O=0.6217639481372091; P=0.6217639481372091;
     Q=692.5187400738176; R=6232.668660664358; S=0.5537744178715619; T=0.5537744178715618 },
  @{ Row=3;  A="FAPs";  B="Angpt1"; C="Tek"; D="FAPs";
     E=3; F=1; G=12.95459633333333; H=38.863789; I=0.8906505749177925; J=0.8906505749177924;
     K=3; L=1; M=31.083557; N=93.25067100000001; O=0.3615335470438062; P=0.3615335470438062;
     Q=402.6749335391577; R=3624.074401852419; S=0.3220000615266349; T=0.3220000615266348 },
  @{ Row=4;  A="FAPs";  B="Angpt1"; C="Tek"; D="MuSCs";
     E=3; F=1; G=12.95459633333333; H=38.863789; I=0.8906505749177925; J=0.8906505749177924;
     K=3; L=1; M=1.436030666666667; N=4.308092; O=0.01670250481898457; P=0.01670250481898457;
     Q=18.60319760895422; R=167.428778480588; S=0.01487609551959581; T=0.01487609551959581 },
  @{ Row=5;  A="MuSCs"; B="Angpt1"; C="Tek"; D="ECs";
     E=3; F=1; G=1.590497666666667; H=4.771493; I=0.1093494250822076; J=0.1093494250822076;
     K=3; L=1; M=53.457377; N=160.372131; O=0.6217639481372091; P=0.6217639481372091;
     Q=85.02383338462035; R=765.214500461583; S=0.06798953026564734; T=0.06798953026564734 },
  @{ Row=6;  A="MuSCs"; B="Angpt1"; C="Tek"; D="FAPs";
     E=3; F=1; G=1.590497666666667; H=4.771493; I=0.1093494250822076; J=0.1093494250822076;
     K=3; L=1; M=31.083557; N=93.25067100000001; O=0.3615335470438062; P=0.3615335470438062;
     Q=49.43832488020034; R=444.9449239218031; S=0.03953348551717146; T=0.03953348551717146 },
  @{ Row=7;  A="MuSCs"; B="Angpt1"; C="Tek"; D="MuSCs";
     E=3; F=1; G=1.590497666666667; H=4.771493; I=0.1093494250822076; J=0.1093494250822076;
     K=3; L=1; M=1.436030666666667; N=4.308092; O=0.01670250481898457; P=0.01670250481898457;
     Q=2.284003424595111; R=20.556030821356; S=0.001826409299388765; T=0.001826409299388765 }
)

foreach ($row in $data) {
  $r = $row.Row
  $ws.Cells.Item($r, 1).Value = $row.A
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  $ws.Cells.Item($r, 6).Value = $row.F
  $ws.Cells.Item($r, 7).Value = $row.G
  $ws.Cells.Item($r, 8).Value = $row.H
  $ws.Cells.Item($r, 9).Value = $row.I
  $ws.Cells.Item($r, 10).Value = $row.J
  $ws.Cells.Item($r, 11).Value = $row.K
  $ws.Cells.Item($r, 12).Value = $row.L
  $ws.Cells.Item($r, 13).Value = $row.M
  $ws.Cells.Item($r, 14).Value = $row.N
  $ws.Cells.Item($r, 15).Value = $row.O
  $ws.Cells.Item($r, 16).Value = $row.P
  $ws.Cells.Item($r, 17).Value = $row.Q
  $ws.Cells.Item($r, 18).Value = $row.R
  $ws.Cells.Item($r, 19).Value = $row.S
  $ws.Cells.Item($r, 20).Value = $row.T
}
